$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2865; $ws.Range("I51").Value = 1437.5; $ws.Range("J51").Value = 3816.6667; $ws.Range("K51").Value = 1437.5; $ws.Range("L51").Value = 3816.6667; $ws.Range("M51").Value = -953.5; $ws.Range("N51").Value = -4784.6667
$ws.Range("H76").Value = 3056.24; $ws.Range("I76").Value = 2910.375; $ws.Range("K76").Value = 2910.375; $ws.Range("M76").Value = -2595.375
$ws.Range("H79").Value = 3056.24; $ws.Range("I79").Value = 2910.375; $ws.Range("K79").Value = 2910.375; $ws.Range("M79").Value = -1818.375
$ws.Range("H87").Value = 80000; $ws.Range("I87").Value = 45000; $ws.Range("J87").Value = 115000; $ws.Range("K87").Value = 45000; $ws.Range("L87").Value = 115000; $ws.Range("M87").Value = -43752; $ws.Range("N87").Value = -117496
$ws.Range("H90").Value = 80000; $ws.Range("I90").Value = 45000; $ws.Range("J90").Value = 115000; $ws.Range("K90").Value = 135000; $ws.Range("L90").Value = 345000; $ws.Range("M90").Value = -128760; $ws.Range("N90").Value = -357480
$ws.Range("H138").Value = 102133.5; $ws.Range("J138").Value = 2985.6667; $ws.Range("L138").Value = 8957.000100000001; $ws.Range("N138").Value = -19237.0001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 609.8333; $ws.Range("I22").Value = 609.8333; $ws.Range("J22").Value = 0; $ws.Range("K22").Value = 609.8333; $ws.Range("L22").Value = 0; $ws.Range("M22").Value = -310.8333
$ws.Range("N22").ClearContents()
$ws.Range("H32").Value = 2877.475; $ws.Range("I32").Value = 2607.7576; $ws.Range("J32").Value = 4149; $ws.Range("K32").Value = 2607.7576; $ws.Range("L32").Value = 4149; $ws.Range("M32").Value = -2320.7576; $ws.Range("N32").Value = -4723
$ws.Range("H61").Value = 79907; $ws.Range("I61").Value = 2977; $ws.Range("K61").Value = 2977; $ws.Range("M61").Value = -2765
$ws.Range("H74").Value = 5285.8; $ws.Range("I74").Value = 2833.5557; $ws.Range("J74").Value = 7292.1816; $ws.Range("K74").Value = 2833.5557; $ws.Range("L74").Value = 7292.1816; $ws.Range("M74").Value = -1959.5557; $ws.Range("N74").Value = -9040.1816
$ws.Range("H77").Value = 5285.8; $ws.Range("I77").Value = 2833.5557; $ws.Range("J77").Value = 7292.1816; $ws.Range("K77").Value = 14167.7785; $ws.Range("L77").Value = 36460.908; $ws.Range("M77").Value = -9799.7785; $ws.Range("N77").Value = -45196.908
$ws.Range("H97").Value = 1168.7; $ws.Range("J97").Value = 2000; $ws.Range("L97").Value = 2000; $ws.Range("N97").Value = -2992
$ws.Range("H102").Value = 63591.777; $ws.Range("I102").Value = 69210.07000000001; $ws.Range("K102").Value = 69210.07000000001; $ws.Range("M102").Value = -67588.07000000001
$ws.Range("H122").Value = 1174.8462; $ws.Range("I122").Value = 1174.8462; $ws.Range("K122").Value = 3524.5386; $ws.Range("M122").Value = -1074.5386
$ws.Range("H123").Value = 62500; $ws.Range("J123").Value = 62500; $ws.Range("L123").Value = 62500; $ws.Range("N123").Value = -72300
$ws.Range("H136").Value = 79907; $ws.Range("I136").Value = 2977; $ws.Range("K136").Value = 8931; $ws.Range("M136").Value = -6381

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 62074; $ws.Range("I20").Value = 73088; $ws.Range("K20").Value = 73088; $ws.Range("M20").Value = -72841
$ws.Range("H54").Value = 1895; $ws.Range("J54").Value = 0; $ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H94").Value = 3173.5; $ws.Range("I94").Value = 3107.6; $ws.Range("K94").Value = 3107.6; $ws.Range("M94").Value = -2656.6
$ws.Range("H99").Value = 168360.17; $ws.Range("I99").Value = 250975; $ws.Range("K99").Value = 250975; $ws.Range("M99").Value = -249477
$ws.Range("H105").Value = 93240.37; $ws.Range("I105").Value = 144685.58; $ws.Range("K105").Value = 144685.58; $ws.Range("M105").Value = -142938.58
$ws.Range("H134").Value = 4077.975; $ws.Range("I134").Value = 2556.6365; $ws.Range("J134").Value = 11250; $ws.Range("K134").Value = 7669.9095; $ws.Range("L134").Value = 33750; $ws.Range("M134").Value = -5134.9095; $ws.Range("N134").Value = -38820

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2763.6; $ws.Range("I31").Value = 2218.4167; $ws.Range("J31").Value = 3127.0557; $ws.Range("K31").Value = 2218.4167; $ws.Range("L31").Value = 3127.0557; $ws.Range("M31").Value = -1923.4167; $ws.Range("N31").Value = -3717.0557
$ws.Range("H34").Value = 2763.6; $ws.Range("I34").Value = 2218.4167; $ws.Range("J34").Value = 3127.0557; $ws.Range("K34").Value = 2218.4167; $ws.Range("L34").Value = 3127.0557; $ws.Range("M34").Value = -2016.4167; $ws.Range("N34").Value = -3531.0557
$ws.Range("H82").Value = 87000; $ws.Range("I82").Value = 0; $ws.Range("J82").Value = 87000; $ws.Range("K82").Value = 0; $ws.Range("L82").Value = 87000; $ws.Range("N82").Value = -87722
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 87000; $ws.Range("I85").Value = 0; $ws.Range("J85").Value = 87000; $ws.Range("K85").Value = 0; $ws.Range("L85").Value = 87000; $ws.Range("N85").Value = -89496
$ws.Range("M85").ClearContents()
$ws.Range("H117").Value = 31449.1; $ws.Range("J117").Value = 31449.1; $ws.Range("L117").Value = 31449.1; $ws.Range("N117").Value = -40627.1

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2651; $ws.Range("I36").Value = 2651; $ws.Range("K36").Value = 7953; $ws.Range("M36").Value = -7784
$ws.Range("H109").Value = 2689.923; $ws.Range("I109").Value = 954; $ws.Range("J109").Value = 3774.875; $ws.Range("K109").Value = 2862; $ws.Range("L109").Value = 11324.625; $ws.Range("M109").Value = -1822; $ws.Range("N109").Value = -13404.625

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 610901.6; $ws.Range("I10").Value = 12751.5; $ws.Range("J10").Value = 1009668.3; $ws.Range("K10").Value = 12751.5; $ws.Range("L10").Value = 1009668.3; $ws.Range("M10").Value = -12582.5; $ws.Range("N10").Value = -1010006.3
$ws.Range("H12").Value = 1712.5555; $ws.Range("I12").Value = 1772.7142; $ws.Range("K12").Value = 1772.7142; $ws.Range("M12").Value = -1632.7142
$ws.Range("H36").Value = 1875; $ws.Range("J36").Value = 1500; $ws.Range("L36").Value = 1500; $ws.Range("N36").Value = -2470
$ws.Range("H40").Value = 16330; $ws.Range("I40").Value = 14000; $ws.Range("J40").Value = 17495; $ws.Range("K40").Value = 14000; $ws.Range("L40").Value = 17495; $ws.Range("M40").Value = -13849; $ws.Range("N40").Value = -17797
$ws.Range("H80").Value = 3327; $ws.Range("I80").Value = 2805; $ws.Range("J80").Value = 3457.5; $ws.Range("K80").Value = 2805; $ws.Range("L80").Value = 3457.5; $ws.Range("M80").Value = -1807; $ws.Range("N80").Value = -5453.5
$ws.Range("H83").Value = 3327; $ws.Range("I83").Value = 2805; $ws.Range("J83").Value = 3457.5; $ws.Range("K83").Value = 14025; $ws.Range("L83").Value = 17287.5; $ws.Range("M83").Value = -9033; $ws.Range("N83").Value = -27271.5
$ws.Range("H97").Value = 3874.889; $ws.Range("I97").Value = 3839.2856; $ws.Range("J97").Value = 3999.5; $ws.Range("K97").Value = 3839.2856; $ws.Range("L97").Value = 3999.5; $ws.Range("M97").Value = -3343.2856; $ws.Range("N97").Value = -4991.5
$ws.Range("H102").Value = 1403.3704; $ws.Range("I102").Value = 1329.8334; $ws.Range("K102").Value = 1329.8334; $ws.Range("M102").Value = 292.1666

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2449.6086; $ws.Range("I46").Value = 1777.6666; $ws.Range("J46").Value = 3182.6365; $ws.Range("K46").Value = 1777.6666; $ws.Range("L46").Value = 3182.6365; $ws.Range("M46").Value = -1589.6666; $ws.Range("N46").Value = -3558.6365
$ws.Range("H100").Value = 4591.364; $ws.Range("I100").Value = 4450.5; $ws.Range("J100").Value = 6000; $ws.Range("K100").Value = 4450.5; $ws.Range("L100").Value = 6000; $ws.Range("M100").Value = -3909.5; $ws.Range("N100").Value = -7082
$ws.Range("H136").Value = 2267.84; $ws.Range("I136").Value = 1690.9412; $ws.Range("J136").Value = 3493.75; $ws.Range("K136").Value = 5072.8236; $ws.Range("L136").Value = 10481.25; $ws.Range("M136").Value = -2522.8236; $ws.Range("N136").Value = -15581.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 17858512; $ws.Range("I100").Value = 23810850; $ws.Range("J100").Value = 1497; $ws.Range("K100").Value = 47621700; $ws.Range("L100").Value = 2994; $ws.Range("M100").Value = -47621159; $ws.Range("N100").Value = -4076
